$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The invoice detail template row (row 15) is re-pointed from the
# "HoaDon2" report-source prefix to "HoaDon1" (QLDatPhong fix), keeping
# each column's underlying field (STT/TenPhong/GiaPhong/SoNgay/TongTien)
# assigned to the same columns B..F.
$ws.Range("B15").Value = "%HoaDon1.STT;insert:copystyles"
$ws.Range("C15").Value = "%HoaDon1.TenPhong"
$ws.Range("D15").Value = "%HoaDon1.GiaPhong"
$ws.Range("E15").Value = "%HoaDon1.SoNgay"
$ws.Range("F15").Value = "%HoaDon1.TongTien"

# Move the sheet's saved selection from G15 to F15.
$ws.Range("F15").Select()

# Re-save the built-in "Normal" cell style so it round-trips under its
# English name instead of the localized "Bình thường" caption.
$wb.Styles.Item(1).Delete()
